# Add the new "element rim" (元素圈) game row (row 11) to the Minigame table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data lives in an Excel Table ("表1" / table1.xml) spanning A3:H10.
# Adding a ListRow grows the table range (and its AutoFilter) by one row,
# matching the A3:H10 -> A3:H11 change in the diff.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()
$rng = $newRow.Range()

$rng.Item(1).Value = 17000008
$rng.Item(2).Value = "元素圈"
$rng.Item(3).Value = 14
$rng.Item(4).Value = 21
$rng.Item(5).Value = 27
$rng.Item(6).Value = 1107
$rng.Item(7).Value = "GameButton8"
$rng.Item(8).Value = "t8"

# Match the updated selection (D11) recorded in the worksheet's sheetView.
$null = $ws.Range("D11").Select()
